$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 44713
$ws.Cells.Item(2, 2).Value = 7.597904761904761
$ws.Cells.Item(2, 3).Value = 7.284552574157715
$ws.Cells.Item(3, 1).Value = 44743
$ws.Cells.Item(3, 2).Value = 7.186949999999999
$ws.Cells.Item(3, 3).Value = 6.071340560913086
$ws.Cells.Item(4, 1).Value = 44774
$ws.Cells.Item(4, 2).Value = 8.779478260869567
$ws.Cells.Item(4, 3).Value = 6.072671413421631
$ws.Cells.Item(5, 1).Value = 44805
$ws.Cells.Item(5, 2).Value = 7.757523809523809
$ws.Cells.Item(5, 3).Value = 7.109158992767334
$ws.Cells.Item(6, 1).Value = 44835
$ws.Cells.Item(6, 2).Value = 6.084904761904762
$ws.Cells.Item(6, 3).Value = 6.213035583496094
$ws.Cells.Item(7, 1).Value = 44866
$ws.Cells.Item(7, 2).Value = 6.429761904761905
$ws.Cells.Item(7, 3).Value = 4.075860023498535
$ws.Cells.Item(8, 1).Value = 44896
$ws.Cells.Item(8, 2).Value = 5.768047619047617
$ws.Cells.Item(8, 3).Value = 5.567882061004639
$ws.Cells.Item(9, 1).Value = 44927
$ws.Cells.Item(9, 2).Value = 3.4228
$ws.Cells.Item(9, 3).Value = 3.801314353942871
$ws.Cells.Item(10, 1).Value = 44958
$ws.Cells.Item(10, 2).Value = 2.437473684210526
$ws.Cells.Item(10, 3).Value = 3.741784334182739
$ws.Cells.Item(11, 1).Value = 44986
$ws.Cells.Item(11, 2).Value = 2.407782608695652
$ws.Cells.Item(11, 3).Value = 2.086360692977905
$ws.Cells.Item(12, 1).Value = 45017
$ws.Cells.Item(12, 2).Value = 2.197263157894737
$ws.Cells.Item(12, 3).Value = 2.105981588363647
$ws.Cells.Item(13, 1).Value = 45047
$ws.Cells.Item(13, 2).Value = 2.299318181818181
$ws.Cells.Item(13, 3).Value = 1.92888069152832
$ws.Cells.Item(14, 1).Value = 45078
$ws.Cells.Item(14, 2).Value = 2.474619047619047
$ws.Cells.Item(14, 3).Value = 2.010726451873779
$ws.Cells.Item(15, 1).Value = 45108
$ws.Cells.Item(15, 2).Value = 2.63655
$ws.Cells.Item(15, 3).Value = 2.294023752212524
$ws.Cells.Item(16, 1).Value = 45139
$ws.Cells.Item(16, 2).Value = 2.645130434782609
$ws.Cells.Item(16, 3).Value = 2.311826705932617
$ws.Cells.Item(17, 1).Value = 45170
$ws.Cells.Item(17, 2).Value = 2.69565
$ws.Cells.Item(17, 3).Value = 2.424980640411377
$ws.Cells.Item(18, 1).Value = 45200
$ws.Cells.Item(18, 2).Value = 3.149181818181818
$ws.Cells.Item(18, 3).Value = 2.446173429489136
$ws.Cells.Item(19, 1).Value = 45231
$ws.Cells.Item(19, 2).Value = 3.055523809523809
$ws.Cells.Item(19, 3).Value = 2.876813411712646
$ws.Cells.Item(20, 1).Value = 45261
$ws.Cells.Item(20, 2).Value = 2.53885
$ws.Cells.Item(20, 3).Value = 2.879498481750488
$ws.Cells.Item(21, 1).Value = 45292
$ws.Cells.Item(21, 2).Value = 2.715
$ws.Cells.Item(21, 3).Value = 2.49435019493103
$ws.Cells.Item(22, 1).Value = 45323
$ws.Cells.Item(22, 2).Value = 1.7955
$ws.Cells.Item(22, 3).Value = 2.515231847763062
$ws.Cells.Item(23, 1).Value = 45352
$ws.Cells.Item(23, 2).Value = 1.7473
$ws.Cells.Item(23, 3).Value = 1.766583204269409
$ws.Cells.Item(24, 1).Value = 45383
$ws.Cells.Item(24, 2).Value = 1.791227272727273
$ws.Cells.Item(24, 3).Value = 1.799374341964722
$ws.Cells.Item(25, 1).Value = 45413
$ws.Cells.Item(25, 2).Value = 2.418
$ws.Cells.Item(25, 3).Value = 1.848229765892029
$ws.Cells.Item(26, 1).Value = 45444
$ws.Cells.Item(26, 2).Value = 2.809578947368421
$ws.Cells.Item(26, 3).Value = 2.015389680862427
$ws.Cells.Item(27, 1).Value = 45474
$ws.Cells.Item(27, 2).Value = 2.208681818181818
$ws.Cells.Item(27, 3).Value = 2.792648553848267
$ws.Cells.Item(28, 1).Value = 45505
$ws.Cells.Item(28, 2).Value = 2.086782608695652
$ws.Cells.Item(28, 3).Value = 2.164881467819214
$ws.Cells.Item(29, 1).Value = 45536
$ws.Cells.Item(29, 2).Value = 2.409250000000001
$ws.Cells.Item(29, 3).Value = 2.064095258712769
$ws.Cells.Item(30, 1).Value = 45566
$ws.Cells.Item(30, 2).Value = 2.576956521739131
$ws.Cells.Item(30, 3).Value = 2.078608751296997
$ws.Cells.Item(31, 1).Value = 45597
$ws.Cells.Item(31, 2).Value = 2.982
$ws.Cells.Item(31, 3).Value = 2.265103101730347
$ws.Cells.Item(32, 1).Value = 45627
$ws.Cells.Item(32, 2).Value = 3.406619047619048
$ws.Cells.Item(32, 3).Value = 2.626777648925781
$ws.Cells.Item(33, 1).Value = 45658
$ws.Cells.Item(33, 2).Value = 3.721380952380952
$ws.Cells.Item(33, 3).Value = 3.071043729782104
$ws.Cells.Item(34, 1).Value = 45689
$ws.Cells.Item(34, 2).Value = 3.740947368421053
$ws.Cells.Item(34, 3).Value = 4.086101055145264
$ws.Cells.Item(35, 1).Value = 45717
$ws.Cells.Item(35, 2).Value = 4.137476190476191
$ws.Cells.Item(35, 3).Value = 3.667343616485596
$ws.Cells.Item(36, 1).Value = 45748
$ws.Cells.Item(36, 2).Value = 3.4
$ws.Cells.Item(36, 3).Value = 3.785993576049805
$ws.Cells.Item(37, 1).Value = 45778
$ws.Cells.Item(37, 2).Value = 3.5
$ws.Cells.Item(37, 3).Value = 3.004053115844727
$ws.Cells.Item(38, 1).Value = 45809
$ws.Cells.Item(38, 2).Value = 3.7
$ws.Cells.Item(38, 3).Value = 2.942317247390747

$ws.Cells.Item(38, 1).NumberFormat = $ws.Cells.Item(37, 1).NumberFormat
